$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.010.37"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.919.67"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.58"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4594"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3814"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07746"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9778"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.68"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").Value = "1.931.69"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.705"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.955"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07019"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.44"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.69"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "29.029.37"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.355"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "2.162.20"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.071"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.73"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.05"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.615"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.75"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.832"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09321"
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8589"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.084"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.020"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05684"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.155"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.156"
$ws.Range("E39").Value = "  +15.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02043"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.409"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5484"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1753"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.375"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002804"
$ws.Range("E45").Value = "  +2.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.186"
$ws.Range("E46").Value = "  +4.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5184"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.25"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06916"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.19"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.757"
$ws.Range("E51").Value = "  -1.53%  "
